$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Replace the old sample rows (A2:C5) with the new, larger data set
# (partida, columna, fila) that now spans rows 2-17.
$data = @(
    @(5,3,5),
    @(5,4,1),
    @(5,5,2),
    @(5,5,3),
    @(6,1,2),
    @(6,2,2),
    @(6,3,5),
    @(6,5,2),
    @(7,2,1),
    @(7,2,3),
    @(7,3,3),
    @(7,4,2),
    @(8,4,1),
    @(8,4,2),
    @(8,5,1),
    @(8,5,4)
)

# Clear the previous data rows first.
$ws.Range("A2:C5").ClearContents()

# The workbook no longer carries the leftover wide/best-fit formatting
# that used to live on column F (it held nothing but a width override).
$ws.Columns.Item(6).Delete()

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# Move the active selection, matching the saved view state.
$ws.Range("L7").Select() | Out-Null
$ws.Range("N8").Select() | Out-Null
